$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 362, pushing existing rows 362-381 down to 363-382.
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new weekly record.
$ws.Cells.Item(362, 1).Value = 4
$ws.Cells.Item(362, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(362, 3).Value = 'Los Lagos'
$ws.Cells.Item(362, 4).Value = 44931
$ws.Cells.Item(362, 5).Value = 10
$ws.Cells.Item(362, 6).Value = 100112040
$ws.Cells.Item(362, 7).Value = 'Cilantro'
$ws.Cells.Item(362, 8).Value = 'Sin especificar'
$ws.Cells.Item(362, 9).Value = 'Primera'
$ws.Cells.Item(362, 10).Value = 40
$ws.Cells.Item(362, 11).Value = 12000
$ws.Cells.Item(362, 12).Value = 12000
$ws.Cells.Item(362, 13).Value = 12000
$ws.Cells.Item(362, 14).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(362, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(362, 16).Value = 6000
$ws.Cells.Item(362, 17).Value = 2
$ws.Cells.Item(362, 18).Value = 'Hortaliza'
